$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated simulation results for rows 2-25 (case with 380 kV done).
# Column G, M, N (all zeros) and column A (index) are unchanged.
# Each array below holds the new values for rows 2..25 of the given column.

$colB = @(0.6537219968967065,0.6139160630076219,0.5896291302994427,0.5797715301364974,0.5781370931603078,0.589496026196997,0.6399653872213094,0.7401291332955964,0.8144157488218582,0.8483557150263152,0.861228296281098,0.8584550680374718,0.8494143488860288,0.8438792602580634,0.8122005714182592,0.7928036814872428,0.7816609527293394,0.7778906201896234,0.7948670862027711,0.8520692842037363,0.889571891705458,0.8695455516369748,0.7939341934524293,0.7129076574818214)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $row = $i + 2
    $ws.Range("B" + $row).Value2 = $colB[$i]
}

$colC = @(0.07532276124410942,0.07312015266796124,0.07175443171390583,0.07119456716172579,0.07110140216217076,0.07174689460251926,0.07456607862538789,0.07998802989306597,0.08390585438475284,0.08567376906771074,0.08634115256619168,0.08619751286798305,0.08572871710336472,0.08544129349607488,0.08379002640539568,0.08277333919741636,0.08218722107601195,0.08198854112980314,0.08288170684960505,0.08586647065638431,0.08780500647947065,0.08677149701620124,0.0828327188496587,0.07853272712524273)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $row = $i + 2
    $ws.Range("C" + $row).Value2 = $colC[$i]
}

$colD = @(0.1200858742501865,0.1180219903983684,0.1167940951046162,0.1163036776228878,0.1162228477250196,0.1167874407455969,0.1193661224794695,0.1247323883475957,0.1288605913721739,0.1307783039195982,0.1315101488163748,0.131352282812216,0.1308384004331629,0.1305243661863216,0.1287360586441366,0.1276491263740951,0.1270276994398927,0.1268179407863386,0.1277644448017838,0.1309891873989812,0.1331296358850693,0.1319842517916072,0.1277122985435284,0.1232478561650936)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $row = $i + 2
    $ws.Range("D" + $row).Value2 = $colD[$i]
}

$colE = @(0.1232652603590871,0.1233192762499531,0.1233928815279572,0.1234330743918814,0.1234403651897704,0.1233933822521305,0.1232755092398623,0.123364025940127,0.1236224458271664,0.1237816927377686,0.1238479627023423,0.1238334253649462,0.1237870253902358,0.1237593802529773,0.1236128751546772,0.123533656703934,0.1234920181412669,0.1234785951646984,0.1235416835443601,0.1238004924816529,0.1240044039315933,0.123892400943646,0.123538042443684,0.1233060116639706)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $row = $i + 2
    $ws.Range("E" + $row).Value2 = $colE[$i]
}

$colF = @(1.86117351845887,1.870935651370026,1.877749120894769,1.880731935628361,1.881239694021872,1.877788512723427,1.864369540681054,1.844548893484308,1.833935473120334,1.829962555899698,1.828580907772974,1.828873011059784,1.829846426536875,1.830458660144735,1.834212306966919,1.836733963027633,1.838264856423415,1.838797022194584,1.836457197813019,1.829557179116065,1.825763368683042,1.827722759906351,1.836582070546918,1.849216724904593)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $row = $i + 2
    $ws.Range("F" + $row).Value2 = $colF[$i]
}

$colH = @(0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429)
for ($i = 0; $i -lt $colH.Length; $i++) {
    $row = $i + 2
    $ws.Range("H" + $row).Value2 = $colH[$i]
}

$colI = @(1.203488368999068,1.21344278487879,1.220060256618247,1.222884101321675,1.223360681745163,1.220097824959225,1.206815799451483,1.184776582286911,1.171022323251794,1.165293556979378,1.163200093235808,1.163647583367634,1.165119805660375,1.166031466806828,1.171407334639039,1.174840485164538,1.176864854010539,1.177558809935867,1.174469875817252,1.164685318868294,1.158732902096993,1.161869358618382,1.174637270668065,1.190310306174045)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $row = $i + 2
    $ws.Range("I" + $row).Value2 = $colI[$i]
}

$colJ = @(0.1600578772566994,0.1610591200374323,0.1617227044184215,0.162005414992338,0.1620531019893807,0.1617264673464103,0.1603929851967898,0.1581645298034076,0.15676175870394,0.1561742746760935,0.1559590734986571,0.1560050980109224,0.1561564244113534,0.1562500620048652,0.1568011698577365,0.1571522161507275,0.1573588964713011,0.1574296942421647,0.1571143533441379,0.1561117791033944,0.1554988848211778,0.1558221287132646,0.1571314559820642,0.1587261273451261)
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $row = $i + 2
    $ws.Range("J" + $row).Value2 = $colJ[$i]
}

$colK = @(0.406884506053018,0.3695605442179613,0.3466825284027664,0.3373700044117243,0.3358243169546711,0.3465568932801091,0.3940074679549355,0.4873445318900451,0.556069723317222,0.5873625494515977,0.5992160067431485,0.596663005905981,0.5883376734134629,0.5832386080513459,0.5540251965164771,0.5361107558922242,0.5258096506863978,0.5223223764583906,0.5380174930982093,0.5907829346507469,0.6252885375755,0.6068706123789127,0.5371554629335549,0.4620662161180462)
for ($i = 0; $i -lt $colK.Length; $i++) {
    $row = $i + 2
    $ws.Range("K" + $row).Value2 = $colK[$i]
}

$colL = @(0.2914885544889216,0.287081584641804,0.2844843249492754,0.2834533730326783,0.2832838462805469,0.2844703098378432,0.2899465623739985,0.3015423353772917,0.3105783719935431,0.3148001339624074,0.316414678852567,0.3160662547606705,0.3149326465726148,0.314240339537136,0.3103046961762459,0.3079186948495476,0.3065568071609022,0.3060974990770404,0.3081716054210517,0.3152651855732529,0.319993602625118,0.3174615529571128,0.3080572338406995,0.2983142331406867)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $row = $i + 2
    $ws.Range("L" + $row).Value2 = $colL[$i]
}

$colO = @(4.917389103029791,4.951555263342343,4.974739993357673,4.984743170800826,4.986437734572917,4.974872651211399,4.92871183132111,4.855684377287929,4.812678099512084,4.795421749046625,4.789218726655605,4.790539914794579,4.794904778511494,4.797621558899181,4.813852254889156,4.824400089573203,4.83068412370514,4.832849100508355,4.823254776923307,4.793613714927716,4.77617418431015,4.785305228566642,4.823771887642977,4.873569239585265)
for ($i = 0; $i -lt $colO.Length; $i++) {
    $row = $i + 2
    $ws.Range("O" + $row).Value2 = $colO[$i]
}
